# Add the "Solution" worksheet (descriptive statistics) after "phone bills"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Defined name "bills" covering the phone-bill amounts on sheet 1
$wb.Names.Add("bills", "='phone bills'!`$A`$2:`$A`$201")

$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws1)
$ws2.Name = "Solution"

# ---- Header row (bold) ----
$ws2.Range("A1").Value = "Statistic"
$ws2.Range("B1").Value = "Value"
$ws2.Range("C1").Value = "Units"
$ws2.Range("D1").Value = "Interpretation"
$ws2.Range("A1:D1").Font.Bold = $true

# ---- Mean ----
$ws2.Range("A2").Value = "Mean"
$ws2.Range("B2").Formula = "=AVERAGE(bills)"
$ws2.Range("C2").Value = "€"
$ws2.Range("D2").Value = "This is the amount that best represent the sample of bills. "

# ---- Median ----
$ws2.Range("A3").Value = "Median"
$ws2.Range("B3").Formula = "=MEDIAN(bills)"
$ws2.Range("C3").Value = "€"
$ws2.Range("D3").Value = "This is the central value of the distribution. Half of the bills have amounts lower than or equal to this amount and the other half have amounts greater than or equal to it. "

# ---- Mode ----
$ws2.Range("A4").Value = "Mode"
$ws2.Range("B4").Formula = "=MODE(bills)"
$ws2.Range("C4").Value = "€"
$ws2.Range("D4").Value = "This is the most common amount of bills, that is, the most frequent. "

# ---- Quartiles ----
$ws2.Range("A5").Value = "Quartiles"

$ws2.Range("A6").Value = 1
$ws2.Range("B6").Formula = "=QUARTILE(bills,A6)"
$ws2.Range("C6").Value = "€"
$ws2.Range("D6").Value = "25% of the bills have amounts lower than or equal to this amount. "

$ws2.Range("A7").Value = 2
$ws2.Range("B7").Formula = "=QUARTILE(bills,A7)"
$ws2.Range("C7").Value = "€"
$ws2.Range("D7").Value = "50% of the bills have amounts lower than or equal to this amount. It's the same as the median. "

$ws2.Range("A8").Value = 3
$ws2.Range("B8").Formula = "=QUARTILE(bills,A8)"
$ws2.Range("C8").Value = "€"
$ws2.Range("D8").Value = "75% of the bills have amounts lower than or equal to this amount. "

# ---- Percentile 65 ----
$ws2.Range("A9").Value = "Percentile 65"
$ws2.Range("B9").Formula = "=PERCENTILE(bills,0.65)"
$ws2.Range("C9").Value = "€"
$ws2.Range("D9").Value = "65% of the bills have amounts lower than or equal to this amount. "

# ---- Variance ----
$ws2.Range("A10").Value = "Variance"
$ws2.Range("B10").Formula = "=VAR.P(bills)"
$ws2.Range("C10").Value = "'€2"
$ws2.Range("C10").NumberFormat = "[$€-2]\ #,##0;[Red]\-[$€-2]\ #,##0"
$ws2.Range("C10").Characters(2, 1).Font.Superscript = $true
$ws2.Range("D10").Value = "Measures the spread with respect to the mean but it has square units and it's difficult to interpret."
$ws2.Rows.Item(10).RowHeight = 14.25

# ---- Std.Deviation ----
$ws2.Range("A11").Value = "Std.Deviation"
$ws2.Range("B11").Formula = "=STDEV.P(bills)"
$ws2.Range("C11").Value = "€"
$ws2.Range("D11").Value = "Measures the average spread with respecto to the mean. Compared to the range of amounts of the bills in the sample, this value is quite high, so there is enough dispersion with respect to the mean. "
$ws2.Range("C11:D11").Font.Name = "Arial"
$ws2.Range("C11:D11").Font.Size = 10

# ---- Coef.Variation ----
$ws2.Range("A12").Value = "Coef.Variation"
$ws2.Range("B12").Formula = "=B11/ABS(B2)"
$ws2.Range("D12").Value = "This statistics also measures the spread with respect to the mean, but it has no units and it's easier to interpret. As this value is greater than 0.5 that means that there is quite dispersión with respect to the mean and therefore the mean is not very representative of the sample. "
$ws2.Range("D12").Font.Name = "Arial"
$ws2.Range("D12").Font.Size = 10

# ---- Coef.Skewness ----
$ws2.Range("A13").Value = "Coef.Skewness"
$ws2.Range("A13").Font.Name = "Arial"
$ws2.Range("A13").Font.Size = 10
$ws2.Range("B13").Formula = "=SKEW(bills)"
$ws2.Range("D13").Value = "As the value is positive the distribution is right-skewed. That means that there are a lot of bills with lower amounts and few bills with huge amounts. "
$ws2.Range("D13").Font.Name = "Arial"
$ws2.Range("D13").Font.Size = 10

# ---- Coef.Kurtosis ----
$ws2.Range("A14").Value = "Coef.Kurtosis"
$ws2.Range("A14").Font.Name = "Arial"
$ws2.Range("A14").Font.Size = 10
$ws2.Range("B14").Formula = "=KURT(bills)"
$ws2.Range("D14").Value = "As the value is negative this means that the distribution is flatter than a bell curve. "
$ws2.Range("D14").Font.Name = "Arial"
$ws2.Range("D14").Font.Size = 10

# ---- Final interpretation note ----
$ws2.Range("D15").Value = "As both, the Coef. of Skewness and the Coef. Of Kurtosis are between -2 and 2, that means that we can assume that the population of bills is normal. "
$ws2.Range("D15").Font.Name = "Arial"
$ws2.Range("D15").Font.Size = 10

# ---- Column widths (best fit look) ----
$ws2.Columns.Item(1).ColumnWidth = 13.140625
$ws2.Columns.Item(2).ColumnWidth = 11.140625
$ws2.Columns.Item(3).ColumnWidth = 4.7109375

# ---- Page setup ----
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# ---- Selection / active sheet bookkeeping ----
$ws2.Range("D15").Select() | Out-Null
$ws1.Activate() | Out-Null
